$d = $word.ActiveDocument

# Grab the full package OOXML (read-only snapshot) so we can compute, per
# paragraph, what its <w:p>...</w:p> markup should look like once the
# direct-formatting <w:contextualSpacing w:val="0"/> toggle is removed from
# its paragraph properties (w:pPr).
$full = $d.Content.WordOpenXML

# Every top-level paragraph in the body serializes as "<w:p ...>...</w:p>"
# or the bare "<w:p>...</w:p>" (attribute-less) form. This pattern only
# matches that opening tag (not <w:pPr>/<w:pBdr>/etc., since those require
# a following identifier character, not '>' or a space) and is non-greedy
# so each match stops at the first following "</w:p>".
$pRegex = [regex]"<w:p[ >][\s\S]*?</w:p>"
$pMatches = $pRegex.Matches($full)

# Namespaces the extracted fragment needs declared on itself, since in the
# source package it inherited them from the ancestor <w:document>.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

$count = $pMatches.Count
for ($i = 0; $i -lt $count; $i++) {
    $paraXml = $pMatches[$i].Value

    # Drop the direct-formatting contextualSpacing toggle from this
    # paragraph's properties.
    $newParaXml = $paraXml.Replace('<w:contextualSpacing w:val="0"/>', '')

    if ($newParaXml -eq $paraXml) {
        # Nothing to change in this paragraph; skip the InsertXML call.
        continue
    }

    # Make the fragment self-contained so InsertXML can parse it on its
    # own, outside the original package context.
    $newParaXml = $newParaXml -replace '^<w:p>', ("<w:p $wNs $w14Ns>")
    $newParaXml = $newParaXml -replace '^<w:p ', ("<w:p $wNs $w14Ns ")

    $para = $d.Paragraphs($i + 1)
    $para.Range.InsertXML($newParaXml)
}
